$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.45
$ws.Range("I3").Value = 6.5
$ws.Range("L3").Value = 1.17
$ws.Range("M3").Value = 5
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = 2.05
$ws.Range("T3").Value = 9
$ws.Range("AA3").Value = 9
$ws.Range("AF3").Value = 41
# Row 4
$ws.Range("AA4").Value = 8
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 29
# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 2.8
$ws.Range("L5").Value = 1.13
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 1.44
$ws.Range("O5").Value = 2.7
$ws.Range("W5").Value = 21
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 21
# Row 6
$ws.Range("G6").Value = 1.75
# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("I7").Value = 3.5
$ws.Range("AI7").Value = 26
# Row 10
$ws.Range("G10").Value = 1.62
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 5
$ws.Range("N10").Value = 1.7
$ws.Range("O10").Value = 2.1
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 8.5
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 23
$ws.Range("Z10").Value = 13
$ws.Range("AD10").Value = 201
$ws.Range("AI10").Value = 41
$ws.Range("AJ10").Value = 41
# Row 11
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 3.75
$ws.Range("I11").Value = 5.5
$ws.Range("N11").Value = 1.93
$ws.Range("O11").Value = 1.88
$ws.Range("W11").Value = 12
# Row 12
$ws.Range("G12").Value = 1.48
$ws.Range("H12").Value = 4.5
$ws.Range("I12").Value = 6.25
$ws.Range("T12").Value = 7.5
$ws.Range("X12").Value = 12
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 8.5
$ws.Range("AB12").Value = 17
$ws.Range("AE12").Value = 17
# Row 14
$ws.Range("J14").Value = 1.05
$ws.Range("K14").Value = 11
$ws.Range("N14").Value = 1.9
$ws.Range("O14").Value = 1.9
# Row 15
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 1.91
$ws.Range("W15").Value = 13
$ws.Range("AI15").Value = 34
# Row 18
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11
# Row 19
$ws.Range("N19").Value = 1.9
$ws.Range("O19").Value = 1.9
# Row 22
$ws.Range("J22").Value = 1.02
$ws.Range("K22").Value = 19
# Row 23
$ws.Range("G23").Value = 2.5
$ws.Range("I23").Value = 2.9
$ws.Range("J23").Value = 1.07
$ws.Range("K23").Value = 9
$ws.Range("T23").Value = 8
$ws.Range("U23").Value = 12
$ws.Range("W23").Value = 23
$ws.Range("X23").Value = 21
$ws.Range("Y23").Value = 29
$ws.Range("AE23").Value = 9
$ws.Range("AF23").Value = 15
$ws.Range("AG23").Value = 11
$ws.Range("AH23").Value = 29
# Row 24
$ws.Range("N24").Value = 1.9
$ws.Range("O24").Value = 1.9
# Row 25
$ws.Range("L25").Value = 1.14
$ws.Range("M25").Value = 5.5
# Row 27
$ws.Range("K27").Value = 17
$ws.Range("L27").Value = 1.17
$ws.Range("M27").Value = 5
$ws.Range("N27").Value = 1.57
$ws.Range("O27").Value = 2.35
